$d = $word.ActiveDocument

# Locate the "Methods" heading paragraph and insert the new content right
# after it (and before "Results"), matching the target diff.
$methodsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Methods") {
        $methodsPara = $p
        break
    }
}

if ($null -eq $methodsPara) {
    throw "Could not locate the Methods heading paragraph"
}

# Create a fresh empty paragraph right after the Methods heading; this is
# where the new body content will live. Using InsertParagraphAfter avoids
# disturbing the Methods heading's own run/bookmark structure.
$methodsPara.Range.InsertParagraphAfter()
$newPara = $methodsPara.Next()
$insertRange = $newPara.Range
# Exclude the trailing paragraph mark so InsertXML doesn't merge with /
# swallow the following (Results) paragraph.
$insertRange.SetRange($insertRange.Start, $insertRange.End - 1)

$bodyXml = @'
<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The current methodology provided by the author for estimating house price supply function uses parametric models. The author in the paper aims to establish a function</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>r</m:t></m:r><m:r><m:t>(</m:t></m:r><m:r><m:t>ν</m:t></m:r><m:r><m:t>)</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">which relates land price</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>(</m:t></m:r><m:sSub><m:e><m:r><m:t>p</m:t></m:r></m:e><m:sub><m:r><m:t>l</m:t></m:r></m:sub></m:sSub><m:r><m:t>)</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and home value per unit land</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>(</m:t></m:r><m:r><m:t>ν</m:t></m:r><m:r><m:t>)</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">. The parametric functions used for the purpose is OLS with different variations in it. The paper has used multiple transformation on data to produce linear, log linear and polynomial linear models to best identify the relation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Since, the function</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>r</m:t></m:r><m:r><m:t>(</m:t></m:r><m:r><m:t>ν</m:t></m:r><m:r><m:t>)</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">forms the root of further analysis and calculation, it is imperative that we replicate the study of these models to test the robustness and the assumptions behind the models.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">We decided two split the method of implementation in two direction. The first method is to use a different model for replication and the other is changing the loss function used in the model.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">For the first case of different model, we decided use Generalized Linear Model with different distribution families. Using different families allows us to verify the condition of normality of error in the models used by authors. We noticed that the gaussian distribution family with log independent and dependent variable produced a line that fitted the best with the data and the corresponding coefficients were very similar to the log linear model used by the author.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">In the case of different loss function, we wrote a function code to implement log linear regression with gradient descent loss function. The gradient descent loss function provides a more flexible approach because of the presence of hyperparameters</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">learning rate</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">number of iteration</w:t></w:r><w:r><w:t xml:space="preserve">. The author of the paper has not explored this methodology, therefore we had no information about the value of hyperparameters</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">learning rate</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">number of iteration</w:t></w:r><w:r><w:t xml:space="preserve">, therefore we ran a simulation study to determine the hyperparameters that delivered values of coefficients very similar to that of the models used by the author</w:t></w:r></w:p>
'@

$xml = @"
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body>
$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

[void]$insertRange.InsertXML($xml)
